# Weekly update for "Hortaliza, Macroferia Regional de Talca - Pepino ensalada"
# Inserts 4 new rows of this week's price reports at the top of the data table
# (rows 778-781), pushing the existing historical rows down by 4
# (old 778-793 -> new 782-797), growing the used range from A1:R793 to A1:R797.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows right above the current first data row of this block (778).
# Inserting repeatedly at the same index pushes each previous insertion down,
# ending up with four fresh rows at 778-781 and everything that used to start
# at 778 now starting at 782.
$ws.Rows.Item(778).Insert()
$ws.Rows.Item(778).Insert()
$ws.Rows.Item(778).Insert()
$ws.Rows.Item(778).Insert()

# New row 778
$ws.Cells.Item(778, 1).Value = 5
$ws.Cells.Item(778, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(778, 3).Value = "Maule"
$ws.Cells.Item(778, 4).Value = 45239
$ws.Cells.Item(778, 5).Value = 7
$ws.Cells.Item(778, 6).Value = 100112043
$ws.Cells.Item(778, 7).Value = "Pepino ensalada"
$ws.Cells.Item(778, 8).Value = "Alaska"
$ws.Cells.Item(778, 9).Value = "Primera"
$ws.Cells.Item(778, 10).Value = 100
$ws.Cells.Item(778, 11).Value = 22000
$ws.Cells.Item(778, 12).Value = 22000
$ws.Cells.Item(778, 13).Value = 22000
$ws.Cells.Item(778, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(778, 15).Value = "Región del Maule"
$ws.Cells.Item(778, 16).Value = 367
$ws.Cells.Item(778, 17).Value = 60
$ws.Cells.Item(778, 18).Value = "Hortaliza"

# New row 779
$ws.Cells.Item(779, 1).Value = 5
$ws.Cells.Item(779, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(779, 3).Value = "Maule"
$ws.Cells.Item(779, 4).Value = 45239
$ws.Cells.Item(779, 5).Value = 7
$ws.Cells.Item(779, 6).Value = 100112043
$ws.Cells.Item(779, 7).Value = "Pepino ensalada"
$ws.Cells.Item(779, 8).Value = "Sin especificar"
$ws.Cells.Item(779, 9).Value = "Primera"
$ws.Cells.Item(779, 10).Value = 300
$ws.Cells.Item(779, 11).Value = 15000
$ws.Cells.Item(779, 12).Value = 15000
$ws.Cells.Item(779, 13).Value = 15000
$ws.Cells.Item(779, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(779, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(779, 16).Value = 250
$ws.Cells.Item(779, 17).Value = 60
$ws.Cells.Item(779, 18).Value = "Hortaliza"

# New row 780
$ws.Cells.Item(780, 1).Value = 5
$ws.Cells.Item(780, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(780, 3).Value = "Maule"
$ws.Cells.Item(780, 4).Value = 45239
$ws.Cells.Item(780, 5).Value = 7
$ws.Cells.Item(780, 6).Value = 100112043
$ws.Cells.Item(780, 7).Value = "Pepino ensalada"
$ws.Cells.Item(780, 8).Value = "Sin especificar"
$ws.Cells.Item(780, 9).Value = "Primera"
$ws.Cells.Item(780, 10).Value = 400
$ws.Cells.Item(780, 11).Value = 17000
$ws.Cells.Item(780, 12).Value = 17000
$ws.Cells.Item(780, 13).Value = 17000
$ws.Cells.Item(780, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(780, 15).Value = "Región del Maule"
$ws.Cells.Item(780, 16).Value = 212
$ws.Cells.Item(780, 17).Value = 80
$ws.Cells.Item(780, 18).Value = "Hortaliza"

# New row 781
$ws.Cells.Item(781, 1).Value = 5
$ws.Cells.Item(781, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(781, 3).Value = "Maule"
$ws.Cells.Item(781, 4).Value = 45239
$ws.Cells.Item(781, 5).Value = 7
$ws.Cells.Item(781, 6).Value = 100112043
$ws.Cells.Item(781, 7).Value = "Pepino ensalada"
$ws.Cells.Item(781, 8).Value = "Sin especificar"
$ws.Cells.Item(781, 9).Value = "Segunda"
$ws.Cells.Item(781, 10).Value = 100
$ws.Cells.Item(781, 11).Value = 14000
$ws.Cells.Item(781, 12).Value = 14000
$ws.Cells.Item(781, 13).Value = 14000
$ws.Cells.Item(781, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(781, 15).Value = "Región del Maule"
$ws.Cells.Item(781, 16).Value = 140
$ws.Cells.Item(781, 17).Value = 100
$ws.Cells.Item(781, 18).Value = "Hortaliza"
